# Adding table format (row/column) to Workbook output.
#
# For every ObjTables header cell in the workbook (the workbook-level
# "!!!ObjTables ..." banner on the first sheet, plus each sheet's
# "!!ObjTables type='Data' id='...' ..." table banner), refresh the
# embedded timestamp and append a tableFormat='row' attribute.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # The workbook banner (only present once, row 1 of the first sheet) and
    # the per-table banner (row 1 of every other sheet, row 2 of the first
    # sheet) always live in column A of the header rows.
    for ($r = 1; $r -le 2; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value

        if ($val -eq $null) { continue }
        if ($val.ToString().StartsWith("!!!ObjTables")) {
            # Workbook-level banner: just refresh the date.
            $newVal = [System.Text.RegularExpressions.Regex]::Replace(
                $val.ToString(), "date='[^']*'", "date='2020-03-09 15:31:17'")
            $cell.Value = $newVal
        }
        elseif ($val.ToString().StartsWith("!!ObjTables")) {
            # Per-table banner: refresh the date and append tableFormat.
            $newVal = [System.Text.RegularExpressions.Regex]::Replace(
                $val.ToString(), "date='[^']*'", "date='2020-03-09 15:31:17'")
            if (-not $newVal.Contains("tableFormat=")) {
                $newVal = $newVal + " tableFormat='row'"
            }
            $cell.Value = $newVal
        }
    }
}
